$wb = $excel.ActiveWorkbook

# --- Step 1: turn the existing "总计" sheet (sheetId=2) into the new "2022-Q1" sheet ---
$q4Sheet = $wb.Worksheets.Item(1)
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Name = "2022-Q1"
$totalSheet.Cells.Clear()

# Header row (copy the style used by the other per-quarter sheet, then set the text)
$q4Sheet.Range("B1:H1").Copy()
$totalSheet.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("B1").Value = "基金代码"
$totalSheet.Range("C1").Value = "基金名称"
$totalSheet.Range("D1").Value = "基金规模"
$totalSheet.Range("E1").Value = "股票总仓位"
$totalSheet.Range("F1").Value = "仓位占比"
$totalSheet.Range("G1").Value = "持有市值(亿元)"
$totalSheet.Range("H1").Value = "仓位排名"

# Data row
$q4Sheet.Range("A2").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("A2").Value = 0

$totalSheet.Range("B2:G2").NumberFormat = "@"
$totalSheet.Range("B2").Value = "003659"
$totalSheet.Range("C2").Value = "山西证券策略精选灵活配置混合"
$totalSheet.Range("D2").Value = "0.31"
$totalSheet.Range("E2").Value = "84.52"
$totalSheet.Range("F2").Value = "3.06"
$totalSheet.Range("G2").Value = "0.0095"
$totalSheet.Range("B2:G2").ClearFormats()
$totalSheet.Range("H2").Value = 7

# --- Step 2: create a brand new "总计" sheet right after "2022-Q1" ---
$newTotal = $wb.Worksheets.Add($null, $totalSheet)
$newTotal.Name = "总计"

$totalSheet.Range("B1:D1").Copy()
$newTotal.Range("B1:D1").PasteSpecial(-4122)
$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"

$totalSheet.Range("A2").Copy()
$newTotal.Range("A2").PasteSpecial(-4122)
$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 1
$newTotal.Range("D2").Value = 0.01

$newTotal.Range("A3").Copy()
$newTotal.Range("A3").PasteSpecial(-4122)
$newTotal.Range("A2").Copy()
$newTotal.Range("A3").PasteSpecial(-4122)
$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 4
$newTotal.Range("D3").Value = 1.14
